$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert two new columns (D and E) so the old "Accuracy" column (C)
#        becomes "Train (%)" and the old "AUC" column (D) shifts to F,
#        leaving room for the new "Validation (%)" / "Test (%)" columns.
$ws.Range("D1:E1").EntireColumn.Insert()

# --- 2. Header row
$ws.Range("C2").Value = "Train (%)"
$ws.Range("D2").Value = "Validation (%)"
$ws.Range("E2").Value = "Test (%)"

# --- 3. Data rows, already written in the final (AUC-descending) order.
$data = @(
    @("Bernoulli Naïve Bayes", 72.930000000000007, 72.87, 69.36, 0.74),
    @("Gaussian Naïve Bayes", 70.38, 70.03, 67.52, 0.73),
    @("Stochastic Gradient Descent", 76.98, 79.83, 73.62, 0.72),
    @("Logistic Regression", 80.8, 81.39, 78.3, 0.71),
    @("Gradient Boosting", 83.3, 79.83, 77.16, 0.69),
    @("Support Vector Machine", 81.819999999999993, 81.25, 77.45, 0.68),
    @("Light Gradient Boosting Machine", 88.34, 79.12, 76.31, 0.68),
    @("Random Forest", 99.75, 79.83, 75.040000000000006, 0.67),
    @("Decision Tree", 99.75, 75.849999999999994, 70.209999999999994, 0.65),
    @("K Nearest Neighbours", 83.51, 74.86, 70.209999999999994, 0.65)
)

$r = 3
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $r = $r + 1
}

# --- 4. Sort the table by AUC (F) descending - records the sortState
#        metadata Excel keeps alongside the AutoFilter.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("F2:F12"), 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("B2:F12"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- 5. Re-apply the AutoFilter over the full new range.
$ws.Range("B2:F12").AutoFilter()

# --- 6. Column widths for the new / resized columns.
$ws.Columns.Item(3).ColumnWidth = 10.140625
$ws.Columns.Item(4).ColumnWidth = 14.5703125
$ws.Columns.Item(5).ColumnWidth = 9.42578125

# --- 7. Fix up the hidden _FilterDatabase defined name for this sheet.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "Individual & Ensemble Model!_FilterDatabase") {
        $n.RefersTo = "='Individual & Ensemble Model'!`$B`$2:`$F`$12"
    }
}

# --- 8. Selection as left by the author.
$ws.Range("D8").Select()
